$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-5 to reflect the new roster ordering
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Dwi Nur Aini"

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Khairun Anwar"

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Rendi Imam Saputra"

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "Ardan Mizanul Khoiri"

# Append the new rows for the remaining groups
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "Wilda Nazwatun Nisa"

$ws.Range("A7").Value = 2
$ws.Range("B7").Value = "Dede Ahmad Fauzan"

$ws.Range("A8").Value = 3
$ws.Range("B8").Value = "Mochammad Wafi Nur Jihan"

$ws.Range("A9").Value = 3
$ws.Range("B9").Value = "Hoirul Sambudi"

$ws.Range("A10").Value = 3
$ws.Range("B10").Value = "Dinda Ayuni"

$ws.Range("A11").Value = 4
$ws.Range("B11").Value = "Ramlan"

$ws.Range("A12").Value = 4
$ws.Range("B12").Value = "Fadli Al Masani"

$ws.Range("A13").Value = 4
$ws.Range("B13").Value = "Agung Prayuda"

# Move the active selection to match the saved cursor position
$ws.Range("H14").Select()
